# Update the "Correspond Handoff Datetime" / "Correspond Handback DateTime"
# timestamps recorded for the handback run (commit: "Generate Report for Handback").
#
# Row 2 of each language sheet (zh-cn / de-de) tracks the 4c115336... file:
#   E2 = Correspond Handoff Datetime
#   H2 = Correspond Handback DateTime

$wb = $excel.ActiveWorkbook

$zh = $wb.Worksheets.Item("zh-cn")
$zh.Range("E2").Value = "2016-03-24 21:28:46"
$zh.Range("H2").Value = "2016-03-24 21:29:18"

$de = $wb.Worksheets.Item("de-de")
$de.Range("E2").Value = "2016-03-24 21:28:51"
$de.Range("H2").Value = "2016-03-24 21:29:26"
